# Cell dimensions were not picometers
# Add the missing |to_pm / |to_nm unit-conversion filters to the
# unit-cell length/volume placeholders in the structure report template.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "cell_length_a"; New = "cell_length_a|to_pm" },
    @{ Old = "cell_length_b"; New = "cell_length_b|to_pm" },
    @{ Old = "cell_length_c"; New = "cell_length_c|to_pm" },
    @{ Old = "cell_volume";   New = "cell_volume|to_nm" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2) | Out-Null
}
